$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp text
$ws.Range("A1").Value = "Datos actualizados a 4 de Mayo de 2020 a las 14:03"

# Row 4: Estados Unidos -> Estados Unidos
$ws.Range("B4").Value = 1188870
$ws.Range("C4").Value = 748
$ws.Range("E4").Value = 941670

# Row 13: Iran -> Iran
$ws.Range("B13").Value = 98647
$ws.Range("C13").Value = 1223
$ws.Range("D13").Value = 79379
$ws.Range("E13").Value = 12991
$ws.Range("F13").Value = 2676
$ws.Range("G13").Value = 74
$ws.Range("H13").Value = 6277

# Row 25: Suecia -> Suecia
$ws.Range("B25").Value = 22721
$ws.Range("C25").Value = 404
$ws.Range("E25").Value = 18947
$ws.Range("F25").Value = 455
$ws.Range("G25").Value = 90
$ws.Range("H25").Value = 2769

# Row 32: Austria -> Catar
$ws.Range("A32").Value = "Catar"
$ws.Range("B32").Value = 16191
$ws.Range("C32").Value = 640
$ws.Range("D32").Value = 1810
$ws.Range("E32").Value = 14369
$ws.Range("F32").Value = 72
$ws.Range("H32").Value = 12

# Row 33: Catar -> Austria
$ws.Range("A33").Value = "Austria"
$ws.Range("B33").Value = 15597
$ws.Range("D33").Value = 13316
$ws.Range("E33").Value = 1681
$ws.Range("F33").Value = 114
$ws.Range("G33").Value = 2
$ws.Range("H33").Value = 600

# Row 54: Finlandia -> Finlandia
$ws.Range("D54").Value = 3500
$ws.Range("E54").Value = 1587
$ws.Range("G54").Value = 10
$ws.Range("H54").Value = 240

# Row 55: Marruecos -> Kuwait
$ws.Range("A55").Value = "Kuwait"
$ws.Range("B55").Value = 5278
$ws.Range("C55").Value = 295
$ws.Range("D55").Value = 1947
$ws.Range("E55").Value = 3291
$ws.Range("F55").Value = 79
$ws.Range("G55").Value = 2
$ws.Range("H55").Value = 40

# Row 56: Kuwait -> Marruecos
$ws.Range("A56").Value = "Marruecos"
$ws.Range("B56").Value = 5000
$ws.Range("C56").Value = 97
$ws.Range("D56").Value = 1565
$ws.Range("E56").Value = 3258
$ws.Range("F56").Value = 1
$ws.Range("G56").Value = 3
$ws.Range("H56").Value = 177

# Row 146: Tayikistan -> Zambia
$ws.Range("A146").Value = "Zambia"
$ws.Range("B146").Value = 137
$ws.Range("C146").Value = 13
$ws.Range("D146").Value = 78
$ws.Range("E146").Value = 56
$ws.Range("F146").Value = 1
$ws.Range("H146").Value = 3

# Row 147: Guayana Francesa -> Tayikistan
$ws.Range("A147").Value = "Tayikistan"
$ws.Range("D147").Value = 0
$ws.Range("E147").Value = 126
$ws.Range("F147").Value = 0
$ws.Range("H147").Value = 2

# Row 148: Togo -> Guayana Francesa
$ws.Range("A148").Value = "Guayana Francesa"
$ws.Range("B148").Value = 128
$ws.Range("D148").Value = 98
$ws.Range("E148").Value = 29
$ws.Range("F148").Value = 2
$ws.Range("H148").Value = 1

# Row 149: Zambia -> Togo
$ws.Range("A149").Value = "Togo"
$ws.Range("D149").Value = 67
$ws.Range("E149").Value = 48
$ws.Range("F149").Value = 0
$ws.Range("H149").Value = 9

# Row 192: Gambia -> San Vicente y las Granadinas
$ws.Range("A192").Value = "San Vicente y las Granadinas"
$ws.Range("C192").Value = 1
$ws.Range("E192").Value = 8
$ws.Range("H192").Value = 0

# Row 193: Santo Tome y Principe -> Gambia
$ws.Range("A193").Value = "Gambia"
$ws.Range("B193").Value = 17
$ws.Range("D193").Value = 9
$ws.Range("E193").Value = 7

# Row 194: San Vicente y las Granadinas -> Santo Tome y Principe
$ws.Range("A194").Value = "Santo Tome y Principe"
$ws.Range("D194").Value = 4
$ws.Range("E194").Value = 11
$ws.Range("H194").Value = 1
